$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 168697
$ws.Range("C4").Value = 159550
$ws.Range("C7").Value = 5.42
$ws.Range("C8").Value = 65.42
